$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: some Price (column D) values are numeric-looking strings (e.g. "1.00", "227.33").
# They must stay stored as text, just like in the original workbook, so a leading
# apostrophe is prepended before assignment - this mirrors how Excel stores a value
# typed with a leading apostrophe (quote-prefixed text) and keeps the cell as text
# instead of letting Excel auto-convert it to a number.

$ws.Range("D2").Value = '92.622.58'
$ws.Range("E2").Value = '  -2.63%  '

$ws.Range("D3").Value = '3.292.28'
$ws.Range("E3").Value = '  -4.52%  '

$ws.Range("D4").Value = "`'" + '1.00'
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = "`'" + '227.33'
$ws.Range("E5").Value = '  -5.25%  '

$ws.Range("D6").Value = "`'" + '606.96'
$ws.Range("E6").Value = '  -5.58%  '

$ws.Range("D7").Value = "`'" + '1.36'
$ws.Range("E7").Value = '  -7.50%  '

$ws.Range("D8").Value = "`'" + '0.377'
$ws.Range("E8").Value = '  -6.33%  '

$ws.Range("D9").Value = "`'" + '1.00'
$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("D10").Value = "`'" + '0.916'
$ws.Range("E10").Value = '  -8.14%  '

$ws.Range("D11").Value = '3.291.97'
$ws.Range("E11").Value = '  -4.45%  '

$ws.Range("D12").Value = "`'" + '41.65'
$ws.Range("E12").Value = '  -0.03%  '

$ws.Range("D13").Value = "`'" + '0.191'
$ws.Range("E13").Value = '  -3.24%  '

$ws.Range("B14").Value = 'Toncoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D14").Value = "`'" + '5.92'
$ws.Range("E14").Value = '  -2.78%  '

$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '92.422.47'
$ws.Range("E15").Value = '  -2.65%  '

$ws.Range("D16").Value = '3.904.74'
$ws.Range("E16").Value = '  -4.58%  '

$ws.Range("D17").Value = "`'" + '0.0000240'
$ws.Range("E17").Value = '  -6.24%  '

$ws.Range("D18").Value = "`'" + '7.98'
$ws.Range("E18").Value = '  -5.46%  '

$ws.Range("D19").Value = '3.289.58'
$ws.Range("E19").Value = '  -4.76%  '

$ws.Range("D20").Value = "`'" + '17.03'
$ws.Range("E20").Value = '  -3.91%  '

$ws.Range("D21").Value = "`'" + '10.64'
$ws.Range("E21").Value = '  -6.82%  '

$ws.Range("D22").Value = "`'" + '3.40'
$ws.Range("E22").Value = '  +7.39%  '

$ws.Range("D23").Value = "`'" + '484.97'
$ws.Range("E23").Value = '  -3.38%  '

$ws.Range("D24").Value = "`'" + '0.438'
$ws.Range("E24").Value = '  -13.27%  '

$ws.Range("D25").Value = "`'" + '0.0000177'
$ws.Range("E25").Value = '  -7.65%  '

$ws.Range("D26").Value = "`'" + '6.01'
$ws.Range("E26").Value = '  -8.55%  '

$ws.Range("D27").Value = "`'" + '88.72'
$ws.Range("E27").Value = '  -6.40%  '

$ws.Range("D28").Value = "`'" + '11.59'
$ws.Range("E28").Value = '  -3.36%  '

$ws.Range("D29").Value = '3.450.02'
$ws.Range("E29").Value = '  -4.98%  '

$ws.Range("D30").Value = "`'" + '1.00'
$ws.Range("E30").Value = '  +0.30%  '

$ws.Range("D31").Value = "`'" + '10.93'
$ws.Range("E31").Value = '  -6.87%  '

$ws.Range("D32").Value = "`'" + '0.137'
$ws.Range("E32").Value = '  +0.87%  '

$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").Value = "`'" + '2.58'
$ws.Range("E33").Value = '  -5.89%  '

$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").Value = "`'" + '1.00'
$ws.Range("E34").Value = '  +0.46%  '

$ws.Range("D35").Value = "`'" + '0.171'
$ws.Range("E35").Value = '  -6.87%  '

$ws.Range("D36").Value = "`'" + '27.86'
$ws.Range("E36").Value = '  -10.17%  '

$ws.Range("D37").Value = "`'" + '0.522'
$ws.Range("E37").Value = '  -7.75%  '

$ws.Range("D38").Value = "`'" + '533.85'
$ws.Range("E38").Value = '  +2.48%  '

$ws.Range("E39").Value = '  -0.02%  '

$ws.Range("D40").Value = "`'" + '7.23'
$ws.Range("E40").Value = '  -5.76%  '

$ws.Range("D41").Value = "`'" + '0.146'
$ws.Range("E41").Value = '  -2.65%  '

$ws.Range("D42").Value = "`'" + '1.34'
$ws.Range("E42").Value = '  -6.75%  '

$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = "`'" + '0.847'
$ws.Range("E43").Value = '  -7.06%  '

$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D44").Value = "`'" + '23.90'
$ws.Range("E44").Value = '  -0.75%  '

$ws.Range("B45").Value = 'ImmutableX'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D45").Value = "`'" + '1.66'
$ws.Range("E45").Value = '  -2.08%  '

$ws.Range("B46").Value = 'MantraDAO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D46").Value = "`'" + '3.57'
$ws.Range("E46").Value = '  +2.38%  '

$ws.Range("D47").Value = "`'" + '0.0405'
$ws.Range("E47").Value = '  -2.60%  '

$ws.Range("D48").Value = "`'" + '5.26'
$ws.Range("E48").Value = '  -6.36%  '

$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = "`'" + '2.07'
$ws.Range("E49").Value = '  -3.57%  '

$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").Value = "`'" + '51.58'
$ws.Range("E50").Value = '  -3.38%  '

$ws.Range("D51").Value = "`'" + '7.78'
$ws.Range("E51").Value = '  -2.96%  '
